# refactor: change to PSA naming
# Rename the sector header abbreviations in row 1 (B1:Q1) to the new PSA naming convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "AFF"
$ws.Range("C1").Value = "MAQ"
$ws.Range("D1").Value = "MFG"
$ws.Range("E1").Value = "ESWW"
$ws.Range("F1").Value = "CNS"
$ws.Range("G1").Value = "TRD"
$ws.Range("H1").Value = "TAS"
$ws.Range("I1").Value = "AFSA"
$ws.Range("J1").Value = "IAC"
$ws.Range("K1").Value = "FIA"
$ws.Range("L1").Value = "REOD"
$ws.Range("O1").Value = "EDUC"
$ws.Range("P1").Value = "HHSW"
$ws.Range("Q1").Value = "OS"

$ws.Range("B1:Q1").Select()
